$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all Price (D) cells that will receive numeric-looking
# string values, so Excel does not auto-convert them to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '69.724.28'
$ws.Range('E2').Value = '  -1.76%  '

$ws.Range('D3').Value = '3.689.59'
$ws.Range('E3').Value = '  -2.55%  '

$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = '614.59'
$ws.Range('E5').Value = '  +0.29%  '

$ws.Range('D6').Value = '177.84'
$ws.Range('E6').Value = '  -1.06%  '

$ws.Range('D7').Value = '3.689.46'
$ws.Range('E7').Value = '  -2.46%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  -2.54%  '

$ws.Range('D10').Value = '0.163'
$ws.Range('E10').Value = '  -1.97%  '

$ws.Range('D11').Value = '6.23'
$ws.Range('E11').Value = '  -2.74%  '

$ws.Range('D12').Value = '0.478'
$ws.Range('E12').Value = '  -4.48%  '

$ws.Range('D13').Value = '39.63'
$ws.Range('E13').Value = '  -2.94%  '

$ws.Range('D14').Value = '0.0000251'
$ws.Range('E14').Value = '  -2.47%  '

$ws.Range('D15').Value = '4.310.80'
$ws.Range('E15').Value = '  -2.41%  '

$ws.Range('D16').Value = '3.692.42'
$ws.Range('E16').Value = '  -2.54%  '

$ws.Range('D17').Value = '69.760.38'
$ws.Range('E17').Value = '  -1.93%  '

$ws.Range('E18').Value = '  -1.99%  '

$ws.Range('D19').Value = '7.50'
$ws.Range('E19').Value = '  -0.85%  '

$ws.Range('D20').Value = '16.30'
$ws.Range('E20').Value = '  -2.57%  '

$ws.Range('D21').Value = '499.75'
$ws.Range('E21').Value = '  -4.46%  '

$ws.Range('D22').Value = '9.12'
$ws.Range('E22').Value = '  -3.43%  '

$ws.Range('D23').Value = '0.710'
$ws.Range('E23').Value = '  -4.92%  '

$ws.Range('D24').Value = '2.52'
$ws.Range('E24').Value = '  +1.22%  '

$ws.Range('D25').Value = '86.05'
$ws.Range('E25').Value = '  -2.74%  '

$ws.Range('D26').Value = '11.35'
$ws.Range('E26').Value = '  +2.77%  '

$ws.Range('D27').Value = '12.90'
$ws.Range('E27').Value = '  -5.02%  '

$ws.Range('E28').Value = '  +4.33%  '

$ws.Range('E29').Value = '  -0.05%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '2.88'
$ws.Range('E30').Value = '  -1.28%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '2.42'
$ws.Range('E31').Value = '  -3.81%  '

$ws.Range('D32').Value = '7.91'
$ws.Range('E32').Value = '  -1.55%  '

$ws.Range('D33').Value = '29.98'
$ws.Range('E33').Value = '  -7.07%  '

$ws.Range('E34').Value = '  -1.70%  '

$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.16%  '

$ws.Range('D36').Value = '1.04'
$ws.Range('E36').Value = '  -1.62%  '

$ws.Range('D37').Value = '6.01'
$ws.Range('E37').Value = '  -2.36%  '

$ws.Range('D38').Value = '0.136'
$ws.Range('E38').Value = '  +3.05%  '

$ws.Range('D39').Value = '0.335'
$ws.Range('E39').Value = '  -1.46%  '

$ws.Range('D40').Value = '2.05'
$ws.Range('E40').Value = '  -8.28%  '

$ws.Range('D41').Value = '49.92'
$ws.Range('E41').Value = '  -2.81%  '

$ws.Range('D42').Value = '44.36'
$ws.Range('E42').Value = '  +0.69%  '

$ws.Range('D43').Value = '426.63'
$ws.Range('E43').Value = '  -0.97%  '

$ws.Range('D44').Value = '2.89'
$ws.Range('E44').Value = '  +3.95%  '

$ws.Range('D45').Value = '8.52'
$ws.Range('E45').Value = '  -3.74%  '

$ws.Range('D46').Value = '2.936.57'
$ws.Range('E46').Value = '  -7.13%  '

$ws.Range('E47').Value = '  -2.56%  '

$ws.Range('D48').Value = '27.19'
$ws.Range('E48').Value = '  -2.29%  '

$ws.Range('E49').Value = '  -0.02%  '

$ws.Range('D50').Value = '136.12'
$ws.Range('E50').Value = '  -3.76%  '

$ws.Range('D51').Value = '2.41'
$ws.Range('E51').Value = '  -2.85%  '

# Restore the Normal style on the Price column so number formatting/style
# metadata is not left modified (matches original General-formatted cells).
$priceRange.Style = "Normal"
